$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three discontinued subjects (delete from bottom row to top to keep indices stable)
$ws.Rows("24").Delete()   # PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS
$ws.Rows("23").Delete()   # PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION
$ws.Rows("9").Delete()    # COMPUTO FLEXIBLE (SOFTCOMPUTING)

# Update predicted enrollment (Cupos_Estimados) values for the remaining subjects
$ws.Range("B2").Value = 71   # ADMINISTRACION DE BASES DE DATOS
$ws.Range("B3").Value = 25   # ADMINISTRACION DE REDES
$ws.Range("B4").Value = 43   # ADMINISTRACION DE SERVIDORES
$ws.Range("B5").Value = 37   # ALGORITMIA
$ws.Range("B6").Value = 48   # ALMACENES DE DATOS (DATA WAREHOUSE)
$ws.Range("B7").Value = 43   # BASES DE DATOS
$ws.Range("B8").Value = 58   # CLASIFICACION INTELIGENTE DE DATOS
$ws.Range("B9").Value = 61   # CONTROL DE PROYECTOS
$ws.Range("B10").Value = 81   # ESTADISTICA Y PROCESOS ESTOCASTICOS
$ws.Range("B11").Value = 7   # ESTRUCTURAS DE DATOS I
$ws.Range("B12").Value = 14   # ESTRUCTURAS DE DATOS II
$ws.Range("B13").Value = 41   # HIPERMEDIA
$ws.Range("B14").Value = 11   # INGENIERIA DE SOFTWARE I
$ws.Range("B15").Value = 50   # INGENIERIA DE SOFTWARE II
$ws.Range("B16").Value = 3   # MATEMATICA DISCRETA
$ws.Range("B17").Value = 114   # METODOS MATEMATICOS I
$ws.Range("B18").Value = 97   # METODOS MATEMATICOS II
$ws.Range("B19").Value = 64   # MINERIA DE DATOS
$ws.Range("B20").Value = 106   # PROGRAMACION
$ws.Range("B21").Value = 290   # PROGRAMACION PARA INTERNET
$ws.Range("B22").Value = 80   # SEGURIDAD DE LA INFORMACION
$ws.Range("B23").Value = 71   # SEMINARIO DE SOLUCION DE PROBLEMAS DE BASES DE DATOS
$ws.Range("B24").Value = 94   # SEMINARIO DE SOLUCION DE PROBLEMAS DE PROGRAMACION
$ws.Range("B25").Value = 76   # SEMINARIO DE SOLUCION DE PROBLEMAS DE ALGORITMIA
$ws.Range("B26").Value = 11   # SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS I
$ws.Range("B27").Value = 34   # SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS II
$ws.Range("B28").Value = 122   # SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS I
$ws.Range("B29").Value = 91   # SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS II
$ws.Range("B30").Value = 113   # SEMINARIO DE SOLUCION DE PROBLEMAS DE INGENIERIA DE SOFTWARE I
$ws.Range("B31").Value = 73   # SEMINARIO DE SOLUCION DE PROBLEMAS DE USO, ADAPTACION, EXPLOTACION DE SISTEMAS OPERATIVOS
$ws.Range("B32").Value = 62   # SEMINARIO DE SOLUCION DE PROBLEMAS DE SISTEMAS BASADOS EN CONOCIMIENTO
$ws.Range("B33").Value = 63   # SISTEMAS BASADOS EN CONOCIMIENTO
$ws.Range("B34").Value = 32   # TEORIA DE LA COMPUTACION
$ws.Range("B35").Value = 70   # USO, ADAPTACION Y EXPLOTACION DE SISTEMAS OPERATIVOS
